$d = $word.ActiveDocument

# Target replacement text (Catalan, with curly apostrophes U+2019).
$newText = "Dates de la campanya Cygnus: 10-19 d" + [char]0x2019 + "agost, 9-18 de setembre, del 8 al 17 d" + [char]0x2019 + "octubre"
$marker = "Dates de la campanya"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    $t = $r.Text
    if ($t.Contains($marker)) {
        $start = $r.Start
        $end = $r.End - 1
        $content = $d.Range($start, $end)
        $content.Delete()
        $ins = $d.Range($start, $start)
        $ins.InsertAfter($newText)
    }
}
